$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104
$ws.Range("B104").Value = 7127374
$ws.Range("F104").Value = 'Central Coast Mariners'
$ws.Range("G104").Value = 'Western Sydney Wanderers'
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 'H'
$ws.Range("K104").Value = 1.909
$ws.Range("M104").Value = 3.6
$ws.Range("N104").Value = 2.15
$ws.Range("O104").Value = 3.6
$ws.Range("P104").Value = 3.25
$ws.Range("Q104").Value = -0.25
$ws.Range("R104").Value = 1.86
$ws.Range("S104").Value = 2.04
$ws.Range("T104").Value = 2.75
$ws.Range("U104").Value = 1.975
$ws.Range("V104").Value = 1.875
$ws.Range("W104").Value = 1.15
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.8600000000000001
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.875

# Row 105
$ws.Range("B105").Value = 7127370
$ws.Range("F105").Value = 'Macarthur FC'
$ws.Range("G105").Value = 'Wellington Phoenix'
$ws.Range("I105").Value = 2
$ws.Range("J105").Value = 'A'
$ws.Range("K105").Value = 2.4
$ws.Range("M105").Value = 2.625
$ws.Range("N105").Value = 2.375
$ws.Range("O105").Value = 3.8
$ws.Range("P105").Value = 2.75
$ws.Range("Q105").Value = 0
$ws.Range("R105").Value = 1.8
$ws.Range("S105").Value = 2.05
$ws.Range("T105").Value = 3
$ws.Range("U105").Value = 1.9
$ws.Range("V105").Value = 1.95
$ws.Range("W105").Value = -1
$ws.Range("Y105").Value = 1.75
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 1.05
$ws.Range("AB105").Value = 0
$ws.Range("AC105").Value = -0

# Row 124
$ws.Range("B124").Value = 7127388
$ws.Range("F124").Value = 'Sydney FC'
$ws.Range("G124").Value = 'Brisbane Roar'
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = 'D'
$ws.Range("K124").Value = 1.5
$ws.Range("L124").Value = 5
$ws.Range("M124").Value = 5
$ws.Range("N124").Value = 1.533
$ws.Range("O124").Value = 5.25
$ws.Range("P124").Value = 5
$ws.Range("Q124").Value = -1
$ws.Range("R124").Value = 1.8
$ws.Range("S124").Value = 2.05
$ws.Range("T124").Value = 3.5
$ws.Range("U124").Value = 1.925
$ws.Range("V124").Value = 1.925
$ws.Range("X124").Value = 4.25
$ws.Range("Y124").Value = -1
$ws.Range("AA124").Value = 1.05
$ws.Range("AB124").Value = -1
$ws.Range("AC124").Value = 0.925

# Row 125
$ws.Range("B125").Value = 7128012
$ws.Range("F125").Value = 'Macarthur FC'
$ws.Range("G125").Value = 'Central Coast Mariners'
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 3
$ws.Range("J125").Value = 'A'
$ws.Range("K125").Value = 2.4
$ws.Range("L125").Value = 3.5
$ws.Range("M125").Value = 2.75
$ws.Range("N125").Value = 3.4
$ws.Range("O125").Value = 3.75
$ws.Range("P125").Value = 2.05
$ws.Range("Q125").Value = 0.25
$ws.Range("R125").Value = 2.025
$ws.Range("S125").Value = 1.825
$ws.Range("T125").Value = 3
$ws.Range("U125").Value = 2.05
$ws.Range("V125").Value = 1.8
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 1.05
$ws.Range("AA125").Value = 0.825
$ws.Range("AB125").Value = 0
$ws.Range("AC125").Value = -0

# Row 126
$ws.Range("H126").Value = 7
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 'H'
$ws.Range("O126").Value = 3.6
$ws.Range("R126").Value = 1.9
$ws.Range("S126").Value = 2
$ws.Range("T126").Value = 3
$ws.Range("U126").Value = 1.925
$ws.Range("V126").Value = 1.925
$ws.Range("W126").Value = 0.8500000000000001
$ws.Range("X126").Value = -1
$ws.Range("Y126").Value = -1
$ws.Range("Z126").Value = 0.8999999999999999
$ws.Range("AA126").Value = -1
$ws.Range("AB126").Value = 0.925
$ws.Range("AC126").Value = -1

# Row 127
$ws.Range("H127").Value = 2
$ws.Range("I127").Value = 2
$ws.Range("J127").Value = 'D'
$ws.Range("N127").Value = 3.6
$ws.Range("P127").Value = 1.909
$ws.Range("Q127").Value = 0.5
$ws.Range("R127").Value = 1.875
$ws.Range("S127").Value = 1.975
$ws.Range("W127").Value = -1
$ws.Range("X127").Value = 3
$ws.Range("Y127").Value = -1
$ws.Range("Z127").Value = 0.875
$ws.Range("AA127").Value = -1
$ws.Range("AB127").Value = 0.95
$ws.Range("AC127").Value = -1

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 1
$ws.Range("J128").Value = 'A'
$ws.Range("N128").Value = 2.2
$ws.Range("O128").Value = 4.2
$ws.Range("R128").Value = 2
$ws.Range("S128").Value = 1.85
$ws.Range("U128").Value = 1.875
$ws.Range("V128").Value = 1.975
$ws.Range("W128").Value = -1
$ws.Range("X128").Value = -1
$ws.Range("Y128").Value = 1.75
$ws.Range("Z128").Value = -1
$ws.Range("AA128").Value = 0.8500000000000001
$ws.Range("AB128").Value = -1
$ws.Range("AC128").Value = 0.9750000000000001

# Row 129
$ws.Range("H129").Value = 2
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = 'H'
$ws.Range("N129").Value = 3.1
$ws.Range("O129").Value = 3.5
$ws.Range("P129").Value = 2.2
$ws.Range("R129").Value = 1.95
$ws.Range("S129").Value = 1.95
$ws.Range("T129").Value = 2.75
$ws.Range("U129").Value = 1.925
$ws.Range("V129").Value = 1.925
$ws.Range("W129").Value = 2.1
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = 0.95
$ws.Range("AA129").Value = -1
$ws.Range("AB129").Value = 0.4625
$ws.Range("AC129").Value = -0.5

# Row 130
$ws.Range("H130").Value = 1
$ws.Range("I130").Value = 2
$ws.Range("J130").Value = 'A'
$ws.Range("N130").Value = 1.833
$ws.Range("O130").Value = 4.333
$ws.Range("P130").Value = 3.75
$ws.Range("R130").Value = 1.825
$ws.Range("S130").Value = 2.025
$ws.Range("U130").Value = 1.975
$ws.Range("V130").Value = 1.875
$ws.Range("W130").Value = -1
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 2.75
$ws.Range("Z130").Value = -1
$ws.Range("AA130").Value = 1.025
$ws.Range("AB130").Value = -1
$ws.Range("AC130").Value = 0.875

# Row 131
$ws.Range("H131").Value = 1
$ws.Range("I131").Value = 2
$ws.Range("J131").Value = 'A'
$ws.Range("N131").Value = 2.55
$ws.Range("O131").Value = 4
$ws.Range("P131").Value = 2.45
$ws.Range("Q131").Value = 0
$ws.Range("R131").Value = 2
$ws.Range("S131").Value = 1.85
$ws.Range("T131").Value = 3.25
$ws.Range("U131").Value = 1.85
$ws.Range("V131").Value = 2
$ws.Range("W131").Value = -1
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = 1.45
$ws.Range("Z131").Value = -1
$ws.Range("AA131").Value = 0.8500000000000001
$ws.Range("AB131").Value = -0.5
$ws.Range("AC131").Value = 0.5
